$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the two custom number formats to columns B (#,"K"/"M") and D (#.000,"K"/"M")
$ws.Range("B2:B12").NumberFormat = "[>999999]#,,""M"";[>999]#,""K"";#"
$ws.Range("D2:D12").NumberFormat = "[>999999]#.000,,""M"";[>999]#.000,""K"";#.000"

# Column B: =A<row> formatted with the "K"/"M" format (no decimals)
$ws.Range("B2").Formula = "=A2"
$ws.Range("B3:B12").Formula = "=A3"

# Column D: =A<row> formatted with the "K"/"M" format (3 decimals)
$ws.Range("D2").Formula = "=A2"
$ws.Range("D3:D12").Formula = "=A3"

# Column C: text values showing what column B's format renders as text
$ws.Range("C2").Value = "'1"
$ws.Range("C3").Value = "'10"
$ws.Range("C4").Value = "'102"
$ws.Range("C5").Value = "'102"
$ws.Range("C6").Value = "'1K"
$ws.Range("C7").Value = "'10K"
$ws.Range("C8").Value = "'102K"
$ws.Range("C9").Value = "'1M"
$ws.Range("C10").Value = "'10M"
$ws.Range("C11").Value = "'102M"
$ws.Range("C12").Value = "'1021M"

# Column E: text values showing what column D's format renders as text
# (rows 4 and 10 are filled in after the rest, matching the authored entry order)
$ws.Range("E2").Value = "'1.020"
$ws.Range("E3").Value = "'10.200"
$ws.Range("E5").Value = "'102.102"
$ws.Range("E6").Value = "'1.021K"
$ws.Range("E7").Value = "'10.210K"
$ws.Range("E8").Value = "'102.102K"
$ws.Range("E9").Value = "'1.021M"
$ws.Range("E11").Value = "'102.102M"
$ws.Range("E12").Value = "'1021.021M"
$ws.Range("E4").Value = "'102.000"
$ws.Range("E10").Value = "'10.210M"

# Column E width (closest the UI grid allows to the authored 13.140625 raw width)
$ws.Columns("E").ColumnWidth = 12.3

# Page setup, matching the authored paper size / orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection like the authored file (F8)
$ws.Range("F8").Select()
